# Weekly cryptos-list refresh (GitHub Actions bot).
# Coin (B), Link (C), Price (D) and Volume(1h) (E) columns are re-synced from the
# coinranking.com feed. Two pairs of adjacent rows also swapped rank order:
#   row 19 <-> row 20 (Uniswap now ranks above BitcoinCash)
#   row 42 <-> row 43 (USDe now ranks above RenderToken)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "72.402.05"
$ws.Range("E2").Value2 = "  +4.16%  "
$ws.Range("D3").Value2 = "2.622.24"
$ws.Range("E3").Value2 = "  +4.48%  "
$ws.Range("D4").Value2 = "'1.00"
$ws.Range("E4").Value2 = "  +0.00%  "
$ws.Range("D5").Value2 = "'603.59"
$ws.Range("E5").Value2 = "  +1.11%  "
$ws.Range("D6").Value2 = "'178.55"
$ws.Range("E6").Value2 = "  +1.49%  "
$ws.Range("D7").Value2 = "'1.00"
$ws.Range("E7").Value2 = "  -0.02%  "
$ws.Range("D8").Value2 = "'0.525"
$ws.Range("E8").Value2 = "  +1.47%  "
$ws.Range("D9").Value2 = "2.619.91"
$ws.Range("E9").Value2 = "  +4.42%  "
$ws.Range("E10").Value2 = "  +8.52%  "
$ws.Range("E11").Value2 = "  +0.88%  "
$ws.Range("D12").Value2 = "'0.352"
$ws.Range("E12").Value2 = "  +2.98%  "
$ws.Range("D13").Value2 = "'5.05"
$ws.Range("E13").Value2 = "  +0.78%  "
$ws.Range("D14").Value2 = "3.138.61"
$ws.Range("E14").Value2 = "  +5.83%  "
$ws.Range("D15").Value2 = "'0.0000187"
$ws.Range("E15").Value2 = "  +5.96%  "
$ws.Range("D16").Value2 = "72.268.47"
$ws.Range("E16").Value2 = "  +4.08%  "
$ws.Range("D17").Value2 = "'26.53"
$ws.Range("E17").Value2 = "  +2.52%  "
$ws.Range("D18").Value2 = "2.624.35"
$ws.Range("E18").Value2 = "  +4.11%  "
$ws.Range("B19").Value2 = "Uniswap"
$ws.Range("C19").Value2 = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value2 = "'8.03"
$ws.Range("E19").Value2 = "  +6.29%  "
$ws.Range("B20").Value2 = "BitcoinCash"
$ws.Range("C20").Value2 = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value2 = "'382.78"
$ws.Range("E20").Value2 = "  +5.34%  "
$ws.Range("E21").Value2 = "  +4.96%  "
$ws.Range("D22").Value2 = "'4.17"
$ws.Range("E22").Value2 = "  +2.97%  "
$ws.Range("D23").Value2 = "'2.00"
$ws.Range("E23").Value2 = "  +18.72%  "
$ws.Range("D24").Value2 = "'72.84"
$ws.Range("E24").Value2 = "  +3.28%  "
$ws.Range("E25").Value2 = "  +0.01%  "
$ws.Range("D26").Value2 = "'4.38"
$ws.Range("E26").Value2 = "  +3.66%  "
$ws.Range("D27").Value2 = "'9.86"
$ws.Range("E27").Value2 = "  +9.10%  "
$ws.Range("D28").Value2 = "2.753.45"
$ws.Range("E28").Value2 = "  +4.84%  "
$ws.Range("D29").Value2 = "'0.999"
$ws.Range("E29").Value2 = "  -0.07%  "
$ws.Range("D30").Value2 = "0.0₃0953"
$ws.Range("E30").Value2 = "  +6.51%  "
$ws.Range("D31").Value2 = "'524.93"
$ws.Range("E31").Value2 = "  +2.44%  "
$ws.Range("D32").Value2 = "'8.03"
$ws.Range("E32").Value2 = "  +3.49%  "
$ws.Range("E33").Value2 = "  +6.77%  "
$ws.Range("E34").Value2 = "  +2.44%  "
$ws.Range("D35").Value2 = "'0.999"
$ws.Range("E35").Value2 = "  -0.09%  "
$ws.Range("D36").Value2 = "'164.18"
$ws.Range("E36").Value2 = "  +1.71%  "
$ws.Range("D37").Value2 = "'19.31"
$ws.Range("E37").Value2 = "  +2.89%  "
$ws.Range("D38").Value2 = "'19.10"
$ws.Range("E38").Value2 = "  +1.14%  "
$ws.Range("E39").Value2 = "  +6.02%  "
$ws.Range("D40").Value2 = "'0.111"
$ws.Range("E40").Value2 = "  -6.47%  "
$ws.Range("D41").Value2 = "'1.82"
$ws.Range("E41").Value2 = "  +5.69%  "
$ws.Range("B42").Value2 = "USDe"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").Value2 = "'1.00"
$ws.Range("E42").Value2 = "  -0.06%  "
$ws.Range("B43").Value2 = "RenderToken"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").Value2 = "'5.04"
$ws.Range("E43").Value2 = "  +4.96%  "
$ws.Range("D44").Value2 = "'2.58"
$ws.Range("E44").Value2 = "  +10.95%  "
$ws.Range("D45").Value2 = "'0.332"
$ws.Range("E45").Value2 = "  +3.64%  "
$ws.Range("D46").Value2 = "'39.48"
$ws.Range("E46").Value2 = "  +1.81%  "
$ws.Range("D47").Value2 = "'150.22"
$ws.Range("E47").Value2 = "  +0.24%  "
$ws.Range("D48").Value2 = "'3.68"
$ws.Range("E48").Value2 = "  +2.89%  "
$ws.Range("D49").Value2 = "'0.542"
$ws.Range("E49").Value2 = "  +5.15%  "
$ws.Range("D50").Value2 = "'1.69"
$ws.Range("E50").Value2 = "  +7.69%  "
$ws.Range("D51").Value2 = "0.0₆0262"
$ws.Range("E51").Value2 = "  +3.14%  "
